# "Further additions to section and assignment classes
#  + early framework or GPA calculation
#  + rough final design for system outline"
#
# The backlog sheet gains an H ("Automate tests" column's section total)
# entry for every section, and the running totals for the sections that
# cover "add classes to each semester", "add course info + GPA calc",
# "add grades/weights" and "predict final grade" are revised downward to
# reflect work still outstanding.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 30-33, 36-38, 41-43 and 46-48 are the blank filler rows inside each
# merged section block; they previously had no column-H cell at all. Excel
# backfills them with the same wrap-text formatting used throughout column
# H/G once the column is populated for the section.
$blankRows = @(30, 31, 32, 33, 36, 37, 38, 41, 42, 43, 46, 47, 48)
foreach ($r in $blankRows) {
    $ws.Cells.Item($r, 8).WrapText = $true
}

# Section: "add classes to each semester" (rows 34-35)
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(35, 8).Value = 0

# Section: "add course information and final grades to calculate GPA" (rows 39-40)
$ws.Cells.Item(39, 8).Value = 8
$ws.Cells.Item(40, 8).Value = 3

# Section: "add grades for each class and specify the weight" (rows 44-45)
$ws.Cells.Item(44, 8).Value = 0
$ws.Cells.Item(45, 8).Value = 0

# Section: "predict final grade based on potential grades" (rows 49-50)
$ws.Cells.Item(49, 8).Value = 11
$ws.Cells.Item(50, 8).Value = 4

# Saved view state: scrolled down into the backlog with J39 selected.
$ws.Range("J39").Select()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
